# Update the "JPDC" worksheet's billing rows (row 2 and row 3) per the
# commit's edits: several descriptive text fields changed, and several
# numeric-looking columns (PIDB, GDCN, Internet Bandwidth, Implementation
# Date) were switched from numbers to text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JPDC")

# --- Row 2 ---
$ws.Range("B2").Value = "nlklkjl"
$ws.Range("C2").Value = "Decommissioned"
$ws.Range("D2").Value = "JPDC1"
$ws.Range("E2").Value = "B"

# Columns I, J, M, P on row 2 become text-stored values instead of numbers.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "12"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "345"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "400"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "45387"

# --- Row 3 ---
$ws.Range("C3").Value = "Decommissioned"
$ws.Range("D3").Value = "JPDC1"
$ws.Range("E3").Value = "B"

# Columns I, J, M, P on row 3 become text-stored values instead of numbers.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "345"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "678"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "500"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "45430"
